$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 12.805
$ws.Range("C6").Value = -13.1
$ws.Range("C7").Value = -13.079
$ws.Range("B8").Value = 6.228000000000001
$ws.Range("E8").Value = 12.438
$ws.Range("A12").Value = -21.482
$ws.Range("B12").Value = 6.806999999999999
$ws.Range("B14").Value = 6.679
$ws.Range("C19").Value = -12.5
$ws.Range("D19").Value = -7.797999999999999
$ws.Range("E19").Value = 12.806
$ws.Range("C21").Value = -13.162
$ws.Range("B22").Value = 6.593000000000001
$ws.Range("C24").Value = -12.638
